$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header labels: columns I/J switch from "MAPE" metrics to "MSE" metrics.
# Shared-string bookkeeping (dropping the now-unused "test (MAPE)"/"train
# (MAPE)" strings, renumbering everything else) happens automatically inside
# the engine when the cell text is reassigned.
$ws.Range("I4").Value = "train (MSE)"
$ws.Range("J4").Value = "test (MSE)"

# --- Per-fold test/train MSE values for rows 5-54 (previously blank).
$data = @(
    @(5, 0.000086948533654646157, 0.00042333451549933239),
    @(6, 0.000059528219722124782, 0.00093316449913139236),
    @(7, 0.000081901266101812095, 0.00032540293463383368),
    @(8, 0.000064759207321783311, 0.00077376428739431757),
    @(9, 0.000073732364772458107, 0.000507232657607337),
    @(10, 0.00008257432659467585, 0.00047721279684502348),
    @(11, 0.000085025625631101402, 0.00030885612997153428),
    @(12, 0.000082158327956365104, 0.00032212746043351219),
    @(13, 0.000085563917677034779, 0.00056667974124099258),
    @(14, 0.000077675925777030195, 0.00062279134941526672),
    @(15, 0.0000851048208371345, 0.0003198312917812094),
    @(16, 0.000079132290479764036, 0.0003800980168802119),
    @(17, 0.000089108185076106351, 0.00075213584947885219),
    @(18, 0.000087304919101244808, 0.00058876905930188973),
    @(19, 0.000081841266169339083, 0.00044599114040169311),
    @(20, 0.000077926918470157106, 0.0004292840957872895),
    @(21, 0.000073149868614087518, 0.00034792833112739981),
    @(22, 0.000079538608471441604, 0.00044672917322765351),
    @(23, 0.000077768459332896187, 0.000376535068809573),
    @(24, 0.000065168861784119499, 0.001063817271652787),
    @(25, 0.000058813712148283392, 0.0012260009833994141),
    @(26, 0.000083099818018459988, 0.0003159042951686708),
    @(27, 0.000079526955407847641, 0.00038530711594354452),
    @(28, 0.000077729543528159235, 0.00027796243188387819),
    @(29, 0.000070305543986441313, 0.00044310374667404299),
    @(30, 0.000065092855783523831, 0.00079037682922212453),
    @(31, 0.000071189477410352406, 0.00043306828051111379),
    @(32, 0.000074833881749245944, 0.00029799906827643558),
    @(33, 0.000082926144321294009, 0.00050521526503828032),
    @(34, 0.000081052768986896752, 0.00041339439979002268),
    @(35, 0.000087367440541888009, 0.0002454964572319719),
    @(36, 0.000050158555937379821, 0.0018085064239096191),
    @(37, 0.00008021347638875162, 0.00047021944419683748),
    @(38, 0.000094666772014777631, 0.00035171964948035319),
    @(39, 0.000077195000717061768, 0.00046671366206977498),
    @(40, 0.000068348457991870994, 0.00055503473438678948),
    @(41, 0.000079585709543738263, 0.00031614404718817619),
    @(42, 0.000082391405102818537, 0.00040718918006139252),
    @(43, 0.000071012248666692493, 0.00050741429697194528),
    @(44, 0.000076171526169691232, 0.00048286261344629108),
    @(45, 0.000079363543803281938, 0.00058644215161087431),
    @(46, 0.000066961515769536724, 0.0006705849340518786),
    @(47, 0.000082568593750855082, 0.00039790390870188708),
    @(48, 0.00006613460407317912, 0.000777747879150208),
    @(49, 0.000075172816614106337, 0.00047831356893084402),
    @(50, 0.000093640030222587661, 0.00037208162156254098),
    @(51, 0.000089236230649986141, 0.00037591960734260129),
    @(52, 0.000064614091119626074, 0.00071378199759458),
    @(53, 0.000057599900900623743, 0.0011061662490035699),
    @(54, 0.000091669860630740941, 0.00046063061610157689)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}

# --- Number formatting for the newly populated I/J columns (rows 5-57,
# including the AVERAGE/STDEV.S summary rows which already held formulas):
# scientific notation, centered - matches the new cellXfs entry (numFmtId 11).
$fmtRange = $ws.Range("I5:J57")
$fmtRange.NumberFormat = "0.00E+00"
$fmtRange.HorizontalAlignment = -4108

# --- View state: scrolled down / selection moved to M57.
$ws.Range("M57").Select()

# --- Page setup (paper size + orientation) added to the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
